$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.222.31"
$ws.Range("E2").Value = "  +3.25%  "
$ws.Range("D3").Value = "'1.908.84"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "'326.58"
$ws.Range("E5").Value = "  +3.53%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D7").Value = "'0.5155"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "'0.4022"
$ws.Range("E8").Value = "  +2.23%  "
$ws.Range("D9").Value = "'0.08476"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "'42.73"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "'1.120"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "'23.37"
$ws.Range("E12").Value = "  +13.47%  "
$ws.Range("D13").Value = "'6.467"
$ws.Range("E13").Value = "  +3.48%  "
$ws.Range("D14").Value = "'1.907.03"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "'7.363"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "'95.04"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "'0.00001114"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").Value = "'0.06667"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("D20").Value = "'18.39"
$ws.Range("E20").Value = "  +2.68%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'5.996"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "'30.228.81"
$ws.Range("E23").Value = "  +3.29%  "
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("D25").Value = "'2.211"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "'2.123.29"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "'21.74"
$ws.Range("E27").Value = "  +3.55%  "
$ws.Range("D28").Value = "'161.55"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").Value = "'2.395"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").Value = "'129.64"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").Value = "'1.096"
$ws.Range("E31").Value = "  +3.63%  "
$ws.Range("D32").Value = "'0.1059"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("D33").Value = "'6.052"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "'3.762"
$ws.Range("E34").Value = "  +3.05%  "
$ws.Range("D35").Value = "'0.02505"
$ws.Range("E35").Value = "  +1.04%  "
$ws.Range("D36").Value = "'0.06583"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "'0.2223"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "'5.258"
$ws.Range("E38").Value = "  +2.72%  "
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("D40").Value = "'11.91"
$ws.Range("E40").Value = "  +5.63%  "
$ws.Range("D41").Value = "'8.808"
$ws.Range("E41").Value = "  -3.60%  "
$ws.Range("D42").Value = "'0.6522"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "'1.230"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").Value = "'0.6124"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").Value = "'13.18"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("D46").Value = "'3.722"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").Value = "'1.245"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("D49").Value = "'124.90"
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("D50").Value = "'1.159"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").Value = "'79.40"
$ws.Range("E51").Value = "  +2.00%  "

# Reset number format/style on D cells that were forced to text via apostrophe,
# so they keep the original (unstyled) appearance instead of a new Text style.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D34").PasteSpecial(-4122)
$ws.Range("D35").PasteSpecial(-4122)
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("D43").PasteSpecial(-4122)
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("D51").PasteSpecial(-4122)
